# Apply updated Num_Inclusions values in column C (analysis by cell, "improve algorithm" pass)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "C3"  = 0
    "C5"  = 3
    "C7"  = 15
    "C8"  = 20
    "C9"  = 20
    "C10" = 6
    "C12" = 8
    "C16" = 2
    "C19" = 11
    "C23" = 11
    "C28" = 11
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
